$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 46053.5
$ws.Range("J109").Value = 46053.5
$ws.Range("L109").Value = 46053.5
$ws.Range("N109").Value = -48827.5
$ws.Range("H117").Value = 99863
$ws.Range("J117").Value = 99863
$ws.Range("L117").Value = 99863
$ws.Range("N117").Value = -109041
$ws.Range("H120").Value = 47496.2
$ws.Range("J120").Value = 47496.2
$ws.Range("L120").Value = 47496.2
$ws.Range("N120").Value = -57172.2
$ws.Range("H123").Value = 65405.453
$ws.Range("J123").Value = 65405.453
$ws.Range("L123").Value = 65405.453
$ws.Range("N123").Value = -75205.45300000001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 418.25
$ws.Range("I2").Value = 418.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 418.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -305.25
$ws.Range("H6").Value = 216666.5
$ws.Range("I6").Value = 243999.8
$ws.Range("K6").Value = 243999.8
$ws.Range("M6").Value = -243826.8
$ws.Range("H7").Value = 98990
$ws.Range("J7").Value = 98990
$ws.Range("L7").Value = 98990
$ws.Range("N7").Value = -99218
$ws.Range("H52").Value = 82696
$ws.Range("J52").Value = 82696
$ws.Range("L52").Value = 82696
$ws.Range("N52").Value = -83332
$ws.Range("H108").Value = 88998
$ws.Range("J108").Value = 88998
$ws.Range("L108").Value = 88998
$ws.Range("N108").Value = -96678
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("H116").Value = 418.25
$ws.Range("I116").Value = 418.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 418.25
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1875.75
$ws.Range("H118").Value = 41442.668
$ws.Range("J118").Value = 41442.668
$ws.Range("L118").Value = 41442.668
$ws.Range("N118").Value = -44756.668
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("H121").Value = 54300
$ws.Range("J121").Value = 54300
$ws.Range("L121").Value = 54300
$ws.Range("N121").Value = -57794
$ws.Range("H122").Value = 3942.1667
$ws.Range("I122").Value = 3997.5
$ws.Range("J122").Value = 3499.5
$ws.Range("K122").Value = 11992.5
$ws.Range("L122").Value = 10498.5
$ws.Range("M122").Value = -9542.5
$ws.Range("N122").Value = -15398.5
$ws.Range("H127").Value = 99420.86
$ws.Range("J127").Value = 99420.86
$ws.Range("L127").Value = 99420.86
$ws.Range("N127").Value = -109340.86
$ws.Range("H132").Value = 1966.5834
$ws.Range("I132").Value = 1372.2222
$ws.Range("K132").Value = 4116.6666
$ws.Range("M132").Value = -1586.6666
$ws.Range("N2").Value = $null
$ws.Range("N113").Value = $null
$ws.Range("N116").Value = $null
$ws.Range("N119").Value = $null

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 418.25
$ws.Range("I3").Value = 418.25
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 418.25
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -304.25
$ws.Range("H110").Value = 94944
$ws.Range("J110").Value = 94944
$ws.Range("L110").Value = 94944
$ws.Range("N110").Value = -103124
$ws.Range("H114").Value = 80830
$ws.Range("J114").Value = 80830
$ws.Range("L114").Value = 80830
$ws.Range("N114").Value = -89508
$ws.Range("H115").Value = 79997.28999999999
$ws.Range("J115").Value = 82996.664
$ws.Range("L115").Value = 82996.664
$ws.Range("N115").Value = -86130.664
$ws.Range("H116").Value = 30371
$ws.Range("J116").Value = 30371
$ws.Range("L116").Value = 30371
$ws.Range("N116").Value = -39549
$ws.Range("H118").Value = 71575.164
$ws.Range("J118").Value = 74781
$ws.Range("L118").Value = 74781
$ws.Range("N118").Value = -78095
$ws.Range("H138").Value = 89996.664
$ws.Range("J138").Value = 89996.664
$ws.Range("L138").Value = 89996.664
$ws.Range("N138").Value = -100276.664
$ws.Range("H140").Value = 43518.68
$ws.Range("J140").Value = 43570.5
$ws.Range("L140").Value = 43570.5
$ws.Range("N140").Value = -53930.5
$ws.Range("N3").Value = $null

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 11625000
$ws.Range("I6").Value = 15166667
$ws.Range("K6").Value = 15166667
$ws.Range("M6").Value = -15166554
$ws.Range("H12").Value = 579
$ws.Range("I12").Value = 499.75
$ws.Range("J12").Value = 737.5
$ws.Range("K12").Value = 499.75
$ws.Range("L12").Value = 737.5
$ws.Range("M12").Value = -329.75
$ws.Range("N12").Value = -1077.5
$ws.Range("H18").Value = 27447
$ws.Range("J18").Value = 27447
$ws.Range("L18").Value = 27447
$ws.Range("N18").Value = -27907
$ws.Range("H108").Value = 55469.453
$ws.Range("J108").Value = 55469.453
$ws.Range("L108").Value = 55469.453
$ws.Range("N108").Value = -63149.453
$ws.Range("H114").Value = 40810
$ws.Range("J114").Value = 40810
$ws.Range("L114").Value = 40810
$ws.Range("N114").Value = -49488
$ws.Range("H121").Value = 28858.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2179.6667
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2179.6667
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 6539.000100000001
$ws.Range("N86").Value = -8911.000100000001
$ws.Range("H89").Value = 2179.6667
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2179.6667
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 19617.0003
$ws.Range("N89").Value = -31473.0003
$ws.Range("H139").Value = 3954
$ws.Range("I139").Value = 2067.3333
$ws.Range("J139").Value = 7350
$ws.Range("K139").Value = 6201.999899999999
$ws.Range("L139").Value = 22050
$ws.Range("M139").Value = -1061.999899999999
$ws.Range("N139").Value = -32330
$ws.Range("M86").Value = $null
$ws.Range("M89").Value = $null

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 59800
$ws.Range("J51").Value = 59800
$ws.Range("L51").Value = 59800
$ws.Range("N51").Value = -60818
$ws.Range("H108").Value = 57765.4
$ws.Range("J108").Value = 57765.4
$ws.Range("L108").Value = 57765.4
$ws.Range("N108").Value = -65445.4
$ws.Range("H109").Value = 35763.152
$ws.Range("J109").Value = 35763.152
$ws.Range("L109").Value = 35763.152
$ws.Range("N109").Value = -37843.152
$ws.Range("H110").Value = 74330.09
$ws.Range("J110").Value = 74330.09
$ws.Range("L110").Value = 74330.09
$ws.Range("N110").Value = -82510.09
$ws.Range("H114").Value = 99990
$ws.Range("J114").Value = 99990
$ws.Range("L114").Value = 99990
$ws.Range("N114").Value = -108668
$ws.Range("H116").Value = 59996.668
$ws.Range("J116").Value = 59996.668
$ws.Range("L116").Value = 59996.668
$ws.Range("N116").Value = -69174.66800000001
$ws.Range("H119").Value = 56970.09
$ws.Range("J119").Value = 56398
$ws.Range("L119").Value = 56398
$ws.Range("N119").Value = -66074
$ws.Range("H132").Value = 6628.3335
$ws.Range("I132").Value = 5681.25
$ws.Range("J132").Value = 7710.7144
$ws.Range("K132").Value = 17043.75
$ws.Range("L132").Value = 23132.1432
$ws.Range("M132").Value = -14513.75
$ws.Range("N132").Value = -28192.1432
$ws.Range("H135").Value = 62588.332
$ws.Range("J135").Value = 62588.332
$ws.Range("L135").Value = 62588.332
$ws.Range("N135").Value = -72728.33199999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2455.8333
$ws.Range("J61").Value = 2960.25
$ws.Range("L61").Value = 2960.25
$ws.Range("N61").Value = -3364.25
$ws.Range("H82").Value = 5499
$ws.Range("H85").Value = 5499
$ws.Range("H113").Value = 2455.8333
$ws.Range("J113").Value = 2960.25
$ws.Range("L113").Value = 2960.25
$ws.Range("N113").Value = -7300.25
$ws.Range("H118").Value = 52654.547
$ws.Range("J118").Value = 53920
$ws.Range("L118").Value = 53920
$ws.Range("N118").Value = -57234
$ws.Range("H121").Value = 46856
$ws.Range("J121").Value = 46856
$ws.Range("L121").Value = 46856
$ws.Range("N121").Value = -50350

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7661.8335
$ws.Range("I62").Value = 8020.875
$ws.Range("J62").Value = 6943.75
$ws.Range("K62").Value = 8020.875
$ws.Range("L62").Value = 6943.75
$ws.Range("M62").Value = -7396.875
$ws.Range("N62").Value = -8191.75
$ws.Range("H65").Value = 7661.8335
$ws.Range("I65").Value = 8020.875
$ws.Range("J65").Value = 6943.75
$ws.Range("K65").Value = 40104.375
$ws.Range("L65").Value = 34718.75
$ws.Range("M65").Value = -36984.375
$ws.Range("N65").Value = -40958.75
$ws.Range("H119").Value = 36250
$ws.Range("J119").Value = 36250
$ws.Range("L119").Value = 36250
$ws.Range("N119").Value = -45926
$ws.Range("H132").Value = 13207.37
$ws.Range("I132").Value = 16780.3
$ws.Range("K132").Value = 50340.89999999999
$ws.Range("M132").Value = -47810.89999999999
$ws.Range("H136").Value = 1535.1666
$ws.Range("I136").Value = 1535.1666
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4605.4998
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2055.4998
$ws.Range("N136").Value = $null
